# Update for 12 April
# Adds a new date column (AD) with the "4/11/20" cumulative death counts,
# following the same layout/format as the existing "AC" (4/10/20) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column AC (the previous last data column) into AD,
# so the new column inherits the same borders/number-format/etc.
# (Only rows 1-54 have data; row 55 must stay untouched.)
$ws.Range("AC1:AC54").Copy() | Out-Null
$ws.Range("AD1:AD54").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# New header for the added date column
$ws.Range("AD1").Value = " 4/11/20"

# New cumulative-deaths values for 4/11/20, one per state row (rows 2-54),
# in the same state order as the existing sheet.
$values = @(93,8,108,25,633,274,494,33,47,446,432,5,8,27,677,330,34,55,94,806,19,206,686,1392,64,93,114,6,17,111,23,2183,20,8627,87,7,247,94,51,501,42,56,80,6,101,267,18,25,130,494,6,137,0)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 30).Value = $v
    $row++
}

# Match the saved selection state from the edit (AD2 selected)
$ws.Range("AD2").Select() | Out-Null
